$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4031794617180017
$ws.Range("B3").Value = 0.219527349232547
$ws.Range("B4").Value = 0.09627809665892988
$ws.Range("B5").Value = 0.07071289249117904
$ws.Range("B6").Value = 0.04878684218129294
$ws.Range("B7").Value = 0.04532657944783455
$ws.Range("B8").Value = 0.04529510356837007
$ws.Range("B9").Value = 0.03672005527873305
$ws.Range("B10").Value = 0.03417361942311181
